# Changes "<a-box ...yellow...></a-box>" to "<a-cone ...blue...></a-cone>"
# in the two a-frame code paragraphs, splitting runs exactly the way the
# target document does (separate <w:r> elements per edited span, all
# sharing the paragraph's original rPr).
#
# Technique notes (discovered empirically against this COM-interop
# engine):
#   - Setting Range.Text = "..." followed by toggling Font.Bold on/off
#     forces the edited span to become its own <w:r> with the run's
#     original rPr preserved (the Bold flip is a no-op value-wise but
#     makes the engine commit the run boundary instead of silently
#     re-flowing/merging it back into neighboring text).
#   - If a Range.Text replacement lands *directly adjacent* to a run
#     boundary that was only "softly" created (via a Font toggle with no
#     text change), the replacement can swallow that neighbor back into
#     one run. So: always perform text-changing replacements first, and
#     only afterwards carve out (via a no-text-change Font toggle) the
#     neighboring unchanged spans that must remain separate runs.

$d = $word.ActiveDocument

function Commit-Range($range) {
    # Forces the range to become (or remain) its own run without
    # altering any visible formatting.
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

function Replace-Span($paraIndex, $relStart, $relEnd, $newText) {
    $base = $d.Paragraphs.Item($paraIndex).Range.Start
    $r = $d.Range($base + $relStart, $base + $relEnd)
    $r.Text = $newText
    Commit-Range $r
}

function Isolate-Text($paraIndex, $findText, $length) {
    # Splits off a not-yet-edited span (found by literal text search
    # within the paragraph) into its own run, without changing its text.
    $base = $d.Paragraphs.Item($paraIndex).Range.Start
    $ptext = $d.Paragraphs.Item($paraIndex).Range.Text
    $rel = $ptext.IndexOf($findText)
    $start = $base + $rel
    $r = $d.Range($start, $start + $length)
    Commit-Range $r
}

function Isolate-Char-After($paraIndex, $findText, $offsetWithin) {
    # Splits off a single character located at $offsetWithin chars after
    # the start of $findText, into its own run.
    $base = $d.Paragraphs.Item($paraIndex).Range.Start
    $ptext = $d.Paragraphs.Item($paraIndex).Range.Text
    $rel = $ptext.IndexOf($findText)
    $pos = $base + $rel + $offsetWithin
    $r = $d.Range($pos, $pos + 1)
    Commit-Range $r
}

function Replace-LastSpan($paraIndex, $findText, $length, $newText) {
    $base = $d.Paragraphs.Item($paraIndex).Range.Start
    $ptext = $d.Paragraphs.Item($paraIndex).Range.Text
    $rel = $ptext.LastIndexOf($findText)
    $start = $base + $rel
    $r = $d.Range($start, $start + $length)
    $r.Text = $newText
    Commit-Range $r
}

# ==========================================================================
# Paragraph 15: "<a-box position='0 1 0' color="yellow"></a-box>"
#            -> "<a-cone position ='0 1 0' color="blue "></a-cone>"
#   (this one gains an extra space before each '=' sign)
# ==========================================================================
$p15 = 15
$base15 = $d.Paragraphs.Item($p15).Range.Start

# "box position" -> "cone position " (extra trailing space merges in for now)
$r = $d.Range($base15 + 15, $base15 + 27)
$r.Text = "cone position "
Commit-Range $r

# split the trailing space off into its own run
Isolate-Char-After $p15 "cone position " 13

# isolate "='0 1 0' " as its own run
Isolate-Text $p15 "='0 1 0' " 9

# "yellow" -> "blue " (extra trailing space merges in for now)
$base15b = $d.Paragraphs.Item($p15).Range.Start
$ptext15b = $d.Paragraphs.Item($p15).Range.Text
$yRel = $ptext15b.IndexOf("yellow")
$yStart = $base15b + $yRel
$r = $d.Range($yStart, $yStart + 6)
$r.Text = "blue "
Commit-Range $r

# split the trailing space off into its own run
Isolate-Char-After $p15 "blue " 4

# closing "box" -> "cone"
Replace-LastSpan $p15 "box" 3 "cone"

# isolate the unchanged '"></a-' tail as its own run
Isolate-Text $p15 '"></a-cone>' 6

# ==========================================================================
# Paragraph 19: "<a-box position='0 1 0' color="yellow"></a-box>"
#            -> "<a-cone position='0 1 0' color="blue"></a-cone>"
# ==========================================================================
$p19 = 19
$base19 = $d.Paragraphs.Item($p19).Range.Start

# "box " -> "cone "
$r = $d.Range($base19 + 15, $base19 + 19)
$r.Text = "cone "
Commit-Range $r

# "yellow" -> "blue"
$base19b = $d.Paragraphs.Item($p19).Range.Start
$ptext19b = $d.Paragraphs.Item($p19).Range.Text
$yRel19 = $ptext19b.IndexOf("yellow")
$yStart19 = $base19b + $yRel19
$r = $d.Range($yStart19, $yStart19 + 6)
$r.Text = "blue"
Commit-Range $r

# closing "box" -> "cone"
Replace-LastSpan $p19 "box" 3 "cone"

# isolate the unchanged "position='0 1 0' " span as its own run
Isolate-Text $p19 "position='0 1 0' " 17

# isolate the unchanged '"></a-' tail as its own run
Isolate-Text $p19 '"></a-cone>' 6

Write-Output "Paragraph 15: $($d.Paragraphs.Item(15).Range.Text)"
Write-Output "Paragraph 19: $($d.Paragraphs.Item(19).Range.Text)"
